$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") — copy H1's formatting (bold,
# centered, bordered header style) onto them so they match the rest of
# the header row.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data columns I ("I0") and J ("IF") for rows 2-14.
$data = @{
    2  = @(1, 4)
    3  = @(1, 5)
    4  = @(1, 5)
    5  = @(1, 5)
    6  = @(1, 6)
    7  = @(1, 5)
    8  = @(1, 4)
    9  = @(2, 3)
    10 = @(3, 4)
    11 = @(1, 3)
    12 = @(6, 8)
    13 = @(1, 3)
    14 = @(5, 6)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 9).Value = $vals[0]
    $ws.Cells.Item($r, 10).Value = $vals[1]
}
